{"js": "// Update the signatory's name and NIP (employee ID) wherever they appear in\n// the document body (the \"Nama\" / \"NIP\" lines occur twice: once in the\n// preamble and once in the closing signature block).\nconst body = context.document.body;\n\nconst nameResults = body.search(\"Nuraina\", { matchCase: true, matchWholeWord: false });\nnameResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < nameResults.items.length; i++) {\n  nameResults.items[i].insertText(\"Chusnul Muawanah, S.T., MM.\", \"Replace\");\n}\nawait context.sync();\n\nconst nipResults = body.search(\"19700417 198903 1004\", { matchCase: true, matchWholeWord: false });\nnipResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < nipResults.items.length; i++) {\n  nipResults.items[i].insertText(\"19800104 200901 2004\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Update the signatory's name and NIP (employee ID) wherever they appear in\n# the document (the \"Nama\" / \"NIP\" lines occur twice: once in the preamble\n# and once in the closing signature block).\n$d = $word.ActiveDocument\n$wdReplaceAll = 2\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Text = \"Nuraina\"\n$find.Replacement.Text = \"Chusnul Muawanah, S.T., MM.\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, $wdReplaceAll)\n\n$find.Text = \"19700417 198903 1004\"\n$find.Replacement.Text = \"19800104 200901 2004\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, $wdReplaceAll)\n"}
